$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Column G: per-segment Area formulas (rows 2-15) ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Column H: total area (row 2) ---
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Columns J/K: summary of totals (row 2) ---
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Selection / active cell, matching the saved view state ---
$ws.Range("J2:K2").Select() | Out-Null
